$wb = $excel.ActiveWorkbook

# --- Sheet "Y" (first sheet) ---
$wsY = $wb.Worksheets.Item("Y")

# Add new data in row 2: B2 = "Maize (home consumed)" (new shared string), C2 = 100000000000000
$wsY.Range("B2").Value = "Maize (home consumed)"
$wsY.Range("C2").Value = 100000000000000

# Clear the numeric values from A3:A6 (leave the cells present but empty)
$wsY.Range("A3").ClearContents()
$wsY.Range("A4").ClearContents()
$wsY.Range("A5").ClearContents()
$wsY.Range("A6").ClearContents()

# Set selection on sheet Y to B3
$wsY.Range("B3").Select() | Out-Null

# Activate sheet Y so it becomes the tab-selected sheet
$wsY.Activate() | Out-Null

# --- Sheet "A" ---
$wsA = $wb.Worksheets.Item("A")
$wsA.Range("F2").Select() | Out-Null

# --- Sheet "VA" ---
$wsVA = $wb.Worksheets.Item("VA")
$wsVA.Range("D10").Select() | Out-Null

# Re-activate sheet Y last, to ensure it is the active/selected tab in the saved file
$wsY.Activate() | Out-Null
